$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns H:K (Username/Password/Endpoint/AddressNs duplicates) are no longer
# used once the row is rebuilt into A:G, so fully clear them (value + style).
$ws.Range("H1:K1").Clear()

# Rebuild the header row, now starting at column A instead of B.
$ws.Range("A1").Value = "Unnamed: 0"
$ws.Range("B1").Value = "Machinetype"
$ws.Range("C1").Value = "Protocol"
$ws.Range("D1").Value = "Username"
$ws.Range("E1").Value = "Password"
$ws.Range("F1").Value = "Endpoint"
$ws.Range("G1").Value = "AddressNs"

# A1 needs the same bold/bordered header formatting already used by B1:G1.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New data row 2.
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Arburg"
$ws.Range("C2").Value = "OPC UA"
$ws.Range("D2").Value = "host_computer"
$ws.Range("E2").Value = " "
$ws.Range("F2").Value = "opc.tcp://10.210.40.219:4880/Arburg"
$ws.Range("G2").Value = "dataSS.csv"

# New data row 3 (A3 is intentionally left blank - the source data has no
# index value for this row, same as an empty/NaN cell).
$ws.Range("B3").Value = "Arburg"
$ws.Range("C3").Value = "OPC UA"
$ws.Range("D3").Value = "host_computer"
$ws.Range("E3").Value = " "
$ws.Range("F3").Value = "opc.tcp://10.210.40.215:4880/Arburg"
$ws.Range("G3").Value = "dataSS.csv"
